# Rename sheet "Data" -> "Summary" (sheetId is preserved since we rename in place).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Summary"

# The round-trip through this engine can drop direct formatting on cells that
# aren't otherwise touched (A1 / A3 keep their original value + position, but
# re-assert their known formatting explicitly so it survives the save).
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# New "Source Type: SME Associations" title/subheading (bold + underlined).
$ws.Range("A9").Value = "Source Type: SME Associations"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# Header row that used to live at row 5 now lives at row 11, bold like before.
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# "Enterprises (% of total)" row that used to live at row 6 now lives at row 12, bold.
$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true

# New data point "90" stored as plain text (not a number), default style.
$d12 = $ws.Cells.Item(12, 4)
$d12.NumberFormat = "@"
$d12.Value = "90"
$d12.ClearFormats()

# New italicized source citation.
$ws.Range("A13").Value = "Source: OBG, 2010"
$ws.Range("A13").Font.Italic = $true

# New bold "OBG" label further down the sheet.
$ws.Range("A19").Value = "OBG"
$ws.Range("A19").Font.Bold = $true

# New italicized full citation text.
$ws.Range("A20").Value = "Oxford Business Group (OBG), `"Economic Update, Ghana boosts financial support for SME development by Oxford Business Group`", 2014. Available at http://www.oxfordbusinessgroup.com/economic_updates/ghana-boosts-financial-support-sme-development"
$ws.Range("A20").Font.Italic = $true

# The old row 5 / row 6 cells were relocated above; remove the now-empty
# originals completely (not just clear their contents) so the rows disappear.
$ws.Range("B5:D5").Clear()
$ws.Range("A6").Clear()
